$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Konfiguracja")

# Update the "nie" -> "tak" values in B3:B6 (shared string used by all four cells)
$ws.Range("B3:B6").Value = "tak"

# Update the active selection to match B3:B6
$ws.Activate()
$ws.Range("B3:B6").Select()
